$d = $word.ActiveDocument

# Remove all paragraphs from the second one through the second-to-last
# one (i.e. everything except the title paragraph and the final, empty
# paragraph), collapsing the whole "references" list down to a single
# empty paragraph that will hold the new hyperlink.
$paraCount = $d.Paragraphs.Count
$deleteStart = $d.Paragraphs(2).Range.Start
$deleteEnd = $d.Paragraphs($paraCount - 1).Range.End
$d.Range($deleteStart, $deleteEnd).Delete()

# Turn the now-empty second paragraph into a hyperlink that points at
# the football article.
$url = "https://myrepublica.nagariknetwork.com/news/with-an-impossible-win-klopp-s-liverpool-write-their-own-history/"
$target = $d.Paragraphs(2).Range
$hyperlink = $d.Hyperlinks.Add($target, $url, "", "", $url)

# Re-create the _GoBack bookmark immediately after the hyperlink, same
# as what Word leaves behind at the last edit position.
$hlEnd = $hyperlink.Range.End
$bookmarkRange = $d.Range($hlEnd - 1, $hlEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
